$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("browser")
$ws2 = $wb.Worksheets.Item("#system")
$ws1.Range("A1:O50").Copy()
$ws2.Range("AZ1:BN50").PasteSpecial(-4104)
